# maj7sat: shuffling around the selection of models available.
# trying to balance the uses of saturation
#
# Adds three new "experiment" rows (div-saturation style re-add, moog
# filter add, and a note about moog filter's compiled size) to the
# running size-tracking log on Sheet2, and moves the active selection
# down to A57.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 65: continuation of the existing delta series (no label yet).
$ws.Range("B65").Value = 20172

# Row 66: "readding the div saturation style, removing sin()"
$ws.Range("A66").Value = "readding the div saturation style, removing sin()"
$ws.Range("B66").Value = 20232

# Row 67: "adding moog filter" with accompanying note in column F
$ws.Range("A67").Value = "adding moog filter"
$ws.Range("B67").Value = 20592
$ws.Range("F67").Value = "so yea moog filter consumes 360 bytes of compressed code. Too much to justify as long as biquad exists"

# Move the saved selection on Sheet2 to A57
$ws.Activate() | Out-Null
$ws.Range("A57").Select() | Out-Null
